# Turn the GitHub repo URL on the "GitHub Link" slide into a real hyperlink.
#
# Target slide: the one titled "GitHub Link" whose body text is the literal
# string "https://github.com/Bramhatejareddy/Steganography.git" (slide 9).
# We locate it defensively (by scanning for the text) rather than hard-coding
# the slide index, then add a hyperlink action to the run so PowerPoint
# writes <a:hlinkClick r:id="rIdN"/> into that run's <a:rPr> and adds the
# corresponding External hyperlink relationship to the slide's .rels part.

$p = $ppt.ActivePresentation

$targetUrl = "https://github.com/Bramhatejareddy/Steganography.git"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)

    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)

        if (-not $shape.HasTextFrame) { continue }

        $textRange = $shape.TextFrame.TextRange
        if ($textRange.Text.Trim() -eq $targetUrl) {
            $action = $textRange.ActionSettings.Item(1)
            $hyperlink = $action.Hyperlink
            $hyperlink.Address = $targetUrl
        }
    }
}
